$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# Overview sheet: row 3 is the 819a2cbc file, columns B (zh-cn) and C (de-de) status
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# zh-cn sheet: row 3 is the 819a2cbc file
$zhcn.Range("B3").Value = $newStatus
$zhcn.Range("G3").Value = "2016-03-02 10:05:23"

# de-de sheet: row 3 is the 819a2cbc file
$dede.Range("B3").Value = $newStatus
$dede.Range("G3").Value = "2016-03-02 10:05:44"
